$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific Price cells to Text format to preserve exact string formatting
# (e.g. trailing zeros) since their new values would otherwise be parsed as numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"

# Apply updated cell values per the source diff
$ws.Range("D2").Value = "36.541.83"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "2.097.54"
$ws.Range("E3").Value = "  +10.74%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "248.41"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("E6").Value = "  -3.14%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "45.22"
$ws.Range("E8").Value = "  +4.02%  "
$ws.Range("D9").Value = "61.30"
$ws.Range("E9").Value = "  +7.72%  "
$ws.Range("D10").Value = "0.366"
$ws.Range("E10").Value = "  +1.74%  "
$ws.Range("E11").Value = "  -3.73%  "
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "14.56"
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("D14").Value = "2.405.53"
$ws.Range("E14").Value = "  +10.97%  "
$ws.Range("D15").Value = "0.839"
$ws.Range("E15").Value = "  +5.77%  "
$ws.Range("D16").Value = "2.093.36"
$ws.Range("E16").Value = "  +10.43%  "
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "36.623.13"
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("D19").Value = "72.56"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").Value = "240.77"
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("D22").Value = "12.84"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("D23").Value = "5.03"
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  -10.64%  "
$ws.Range("D26").Value = "169.87"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").Value = "20.53"
$ws.Range("E27").Value = "  +11.12%  "
$ws.Range("D28").Value = "8.92"
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("D29").Value = "1.99"
$ws.Range("E29").Value = "  -8.22%  "
$ws.Range("E30").Value = "  -4.75%  "
$ws.Range("D31").Value = "22.37"
$ws.Range("E31").Value = "  +58.84%  "
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").Value = "0.0910"
$ws.Range("E34").Value = "  +15.58%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "1.88"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "2.31"
$ws.Range("E37").Value = "  +18.42%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "4.06"
$ws.Range("E38").Value = "  -5.40%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "0.905"
$ws.Range("E39").Value = "  +4.78%  "
$ws.Range("D40").Value = "1.35"
$ws.Range("E40").Value = "  -8.90%  "
$ws.Range("E41").Value = "  +8.91%  "
$ws.Range("D42").Value = "99.15"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("E43").Value = "  -4.75%  "
$ws.Range("E44").Value = "  +16.29%  "
$ws.Range("D45").Value = "16.17"
$ws.Range("E45").Value = "  -4.77%  "
$ws.Range("D46").Value = "1.360.68"
$ws.Range("E46").Value = "  +3.32%  "
$ws.Range("E47").Value = "  +2.97%  "
$ws.Range("D48").Value = "2.293.71"
$ws.Range("E48").Value = "  +10.88%  "
$ws.Range("E49").Value = "  +2.49%  "
$ws.Range("E50").Value = "  -3.65%  "
$ws.Range("E51").Value = "  +16.59%  "
